# Germany Verbandsliga base update (17-02-2024 11:11)
# The source feed re-sorted matches that share the same Date within a few
# row-groups by their match id (column B), ascending. Column A (the running
# "id" index) is left untouched; only the data columns B:AC move between rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($row) {
    return $ws.Range("B" + $row + ":AC" + $row).Value()
}

function Set-RowData($row, $data) {
    $ws.Range("B" + $row + ":AC" + $row).Value = $data
}

# --- Rows 38, 39, 40: cyclic rotation -------------------------------------
# after(38) = before(39), after(39) = before(40), after(40) = before(38)
$row38 = Get-RowData 38
$row39 = Get-RowData 39
$row40 = Get-RowData 40

Set-RowData 38 $row39
Set-RowData 39 $row40
Set-RowData 40 $row38

# --- Rows 43, 44: swap -----------------------------------------------------
$row43 = Get-RowData 43
$row44 = Get-RowData 44

Set-RowData 43 $row44
Set-RowData 44 $row43

# --- Rows 111, 112: swap ----------------------------------------------------
$row111 = Get-RowData 111
$row112 = Get-RowData 112

Set-RowData 111 $row112
Set-RowData 112 $row111

# --- Rows 151, 152: swap ----------------------------------------------------
$row151 = Get-RowData 151
$row152 = Get-RowData 152

Set-RowData 151 $row152
Set-RowData 152 $row151
